$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Helper: assign a literal text value to a cell without letting Excel's
# automatic input parsing reinterpret percentage-looking strings (e.g.
# "28.1%") as numeric percentage values. We temporarily place a text
# formula that evaluates to the desired string, then use Copy/PasteSpecial
# (values only) to collapse the formula result back down to a plain text
# value in place, preserving the cell's existing style.
function Set-LiteralText($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# Row 2: reorder "Recorded By" list for the ANATOMY session 1 record
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg"

# Row 3: reorder "Recorded By" list and update attendance count for ANATOMY session 2
$ws.Range("G3").Value = "Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("H3").Value = "103/251"

# Row 10: updated average attendance % statistic
Set-LiteralText $ws.Range("L10") "28.1%"

# Row 15: reorder "Recorded By" list for PARASITOLOGY session 2
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 15: updated average attendance % statistic
Set-LiteralText $ws.Range("S15") "28.1%"
